$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45810.01041666666, 1),
    @(45810.02083333334, 1),
    @(45810.03125, 1),
    @(45810.04166666666, 1),
    @(45810.05208333334, 1),
    @(45810.0625, 1),
    @(45810.07291666666, 1),
    @(45810.08333333334, 1),
    @(45810.09375, 1),
    @(45810.10416666666, 1),
    @(45810.11458333334, 1),
    @(45810.125, 1),
    @(45810.13541666666, 1),
    @(45810.14583333334, 1),
    @(45810.15625, 1),
    @(45810.16666666666, 8),
    @(45810.17708333334, 16),
    @(45810.1875, 18),
    @(45810.19791666666, 18),
    @(45810.20833333334, 23),
    @(45810.21875, 101),
    @(45810.22916666666, 109),
    @(45810.23958333334, 122),
    @(45810.25, 139),
    @(45810.26041666666, 430),
    @(45810.27083333334, 455),
    @(45810.28125, 486),
    @(45810.29166666666, 538),
    @(45810.30208333334, 1029),
    @(45810.3125, 1071),
    @(45810.32291666666, 1121),
    @(45810.33333333334, 1172),
    @(45810.34375, 1587),
    @(45810.35416666666, 1629),
    @(45810.36458333334, 1674),
    @(45810.375, 1723),
    @(45810.38541666666, 1947),
    @(45810.39583333334, 1980),
    @(45810.40625, 2013),
    @(45810.41666666666, 2036),
    @(45810.42708333334, 2190),
    @(45810.4375, 2205),
    @(45810.44791666666, 2218),
    @(45810.45833333334, 2229),
    @(45810.46875, 2283),
    @(45810.47916666666, 2290),
    @(45810.48958333334, 2288),
    @(45810.5, 2287),
    @(45810.51041666666, 2272),
    @(45810.52083333334, 2268),
    @(45810.53125, 2262),
    @(45810.54166666666, 2256),
    @(45810.55208333334, 2203),
    @(45810.5625, 2194),
    @(45810.57291666666, 2182),
    @(45810.58333333334, 2165),
    @(45810.59375, 2064),
    @(45810.60416666666, 2050),
    @(45810.61458333334, 2028),
    @(45810.625, 2003),
    @(45810.63541666666, 1770),
    @(45810.64583333334, 1737),
    @(45810.65625, 1705),
    @(45810.66666666666, 1673),
    @(45810.67708333334, 1349),
    @(45810.6875, 1306),
    @(45810.69791666666, 1266),
    @(45810.70833333334, 1230),
    @(45810.71875, 786),
    @(45810.72916666666, 751),
    @(45810.73958333334, 714),
    @(45810.75, 684),
    @(45810.76041666666, 307),
    @(45810.77083333334, 275),
    @(45810.78125, 255),
    @(45810.79166666666, 237),
    @(45810.80208333334, 60),
    @(45810.8125, 50),
    @(45810.82291666666, 44),
    @(45810.83333333334, 41),
    @(45810.84375, 25),
    @(45810.85416666666, 25),
    @(45810.86458333334, 25),
    @(45810.875, 25),
    @(45810.88541666666, 11),
    @(45810.89583333334, 10),
    @(45810.90625, 10),
    @(45810.91666666666, 10),
    @(45810.92708333334, 1),
    @(45810.9375, 1),
    @(45810.94791666666, 1),
    @(45810.95833333334, 1),
    @(45810.96875, 0),
    @(45810.97916666666, 0),
    @(45810.98958333334, 0),
    @(45811, 0),
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
